# Update "Hjemme passive" values: meanEMG legmaxROM columns (B-E) on rows 1-3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header / ROM values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 - CON data
$ws.Range("B2").Value = 0.17975631303613016
$ws.Range("C2").Value = 0.029867720791489222
$ws.Range("D2").Value = 0.10359765777763366
$ws.Range("E2").Value = 0.009578336710277408

# Row 3 - STR data
$ws.Range("B3").Value = 0.12772286123046353
$ws.Range("C3").Value = 0.04504508733259733
$ws.Range("D3").Value = 0.230627152085778
$ws.Range("E3").Value = 0.030193810726400653
